$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 743
$ws.Range("F3").Value = 72
$ws.Range("F5").Value = 6117
$ws.Range("F7").Value = 195
$ws.Range("F8").Value = 3424
$ws.Range("F9").Value = 700
$ws.Range("F11").Value = 1433
$ws.Range("F12").Value = 4774
$ws.Range("F13").Value = 1781
$ws.Range("F14").Value = 23
$ws.Range("F15").Value = 66
$ws.Range("F17").Value = 217
$ws.Range("F18").Value = 163
$ws.Range("F19").Value = 1062
$ws.Range("F20").Value = 329
$ws.Range("F22").Value = 35
$ws.Range("F24").Value = 5
$ws.Range("F25").Value = 221
$ws.Range("F26").Value = 111
$ws.Range("F27").Value = 21
$ws.Range("F28").Value = 1150
$ws.Range("F29").Value = 432
$ws.Range("F30").Value = 128
$ws.Range("F31").Value = 234
$ws.Range("F33").Value = 28
$ws.Range("F34").Value = 1848
$ws.Range("F35").Value = 2305
$ws.Range("F37").Value = 42
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 301
$ws.Range("F40").Value = 72
$ws.Range("F41").Value = 671
$ws.Range("F42").Value = 470
$ws.Range("F43").Value = 63
$ws.Range("F44").Value = 694
$ws.Range("F45").Value = 51
$ws.Range("F46").Value = 473
$ws.Range("F47").Value = 487
$ws.Range("F48").Value = 240

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 39
$ws.Range("F24").Value = 8

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 827

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 827
$ws.Range("F3").Value = 743
$ws.Range("F4").Value = 72
$ws.Range("F6").Value = 6117
$ws.Range("F8").Value = 195
$ws.Range("F9").Value = 3424
$ws.Range("F11").Value = 1433
$ws.Range("F12").Value = 4774
$ws.Range("F13").Value = 1781
$ws.Range("F14").Value = 23
$ws.Range("F16").Value = 66
$ws.Range("F19").Value = 217
$ws.Range("F20").Value = 163
$ws.Range("F22").Value = 1062
$ws.Range("F23").Value = 329
$ws.Range("F25").Value = 221
$ws.Range("F27").Value = 21
$ws.Range("F28").Value = 1150
$ws.Range("F29").Value = 432
$ws.Range("F30").Value = 128
$ws.Range("F31").Value = 234
$ws.Range("F34").Value = 1848
$ws.Range("F35").Value = 2306
$ws.Range("F37").Value = 20
$ws.Range("F38").Value = 301
$ws.Range("F39").Value = 72
$ws.Range("F42").Value = 671
$ws.Range("F43").Value = 470
$ws.Range("F44").Value = 694
$ws.Range("F45").Value = 473
$ws.Range("F46").Value = 487
$ws.Range("F48").Value = 160

